# Rename the excel sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Misc Exception"
